$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.111.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.01%  "
$ws.Range("D3").Value = "'2.232.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'319.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'101.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.11%  "
$ws.Range("E7").Value = "  -6.98%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.563"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -8.28%  "
$ws.Range("D10").Value = "'37.07"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -9.20%  "
$ws.Range("D11").Value = "'54.45"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.0826"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -9.97%  "
$ws.Range("E13").Value = "  -9.45%  "
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "'2.574.01"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "'0.861"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -12.43%  "
$ws.Range("D17").Value = "'14.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.53%  "
$ws.Range("D18").Value = "'2.231.35"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.76%  "
$ws.Range("D19").Value = "'43.053.05"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -5.10%  "
$ws.Range("D20").Value = "'14.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").Value = "'0.0₃0965"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -9.22%  "
$ws.Range("D22").Value = "'6.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -10.62%  "
$ws.Range("D23").Value = "'65.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -10.99%  "
$ws.Range("D24").Value = "'3.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -10.90%  "
$ws.Range("D25").Value = "'238.01"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -8.75%  "
$ws.Range("D26").Value = "'2.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -8.87%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "'4.07"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "'10.04"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.63%  "
$ws.Range("D30").Value = "'2.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "'6.39"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -13.28%  "
$ws.Range("D32").Value = "'35.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").Value = "'20.48"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.07%  "
$ws.Range("D34").Value = "'0.0874"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.78%  "
$ws.Range("D35").Value = "'153.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.89%  "
$ws.Range("E36").Value = "  -5.05%  "
$ws.Range("D37").Value = "'3.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.74%  "
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  -6.84%  "
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("D41").Value = "'0.103"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -11.75%  "
$ws.Range("D42").Value = "'3.70"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.03%  "
$ws.Range("D43").Value = "'0.0325"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.20%  "
$ws.Range("D44").Value = "'12.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'1.805.82"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "'87.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.93%  "
$ws.Range("E48").Value = "  -9.10%  "
$ws.Range("D49").Value = "'76.78"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.86%  "
$ws.Range("D50").Value = "'5.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("D51").Value = "'59.33"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -15.57%  "
